# "Fixed biomethane in EB's Expected Porfolio"
#
# The "feedstock_to_commodity" sheet had a row (row 8: Biomethane ->
# "Anaerobic Digestion + Conditioning + Gas Engine" -> Electricity) whose
# total (column Z) was driven off a stale external workbook reference
# ([1]portfolio_input!C14*1000) instead of this workbook's Sheet1 inputs.
# The fix removes that row outright (shifting every following row up by
# one), which also drops the now-unused "Anaerobic Digestion +
# Conditioning + Gas Engine" shared string and the now-unused external
# reference/link entirely.

$wb = $excel.ActiveWorkbook

$wsFeedstock = $wb.Worksheets.Item("feedstock_to_commodity")
$wsFeedstock.Activate()
$wsFeedstock.Rows("8:8").Delete()

# Drop the external workbook link that row used to depend on - nothing
# in the workbook references it any more after the row is gone.
$linkSources = $wb.LinkSources(1)
if ($linkSources -ne $null) {
    foreach ($link in $linkSources) {
        $wb.BreakLink($link, 1)
    }
}

# Restore the view state recorded in the saved file: cursor on Sheet1
# moved to A26, and on feedstock_to_commodity moved to the new last row
# (A17, since the sheet now only goes down to row 17).
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Activate()
$wsSheet1.Range("A26").Select()

$wsFeedstock.Activate()
$wsFeedstock.Range("A17").Select()
